$wb = $excel.ActiveWorkbook

# Both sheets get cleaned up to contain the same two FreeText test-question rows
# (replacing the previous three rows of unrelated test data), then the third
# data row is cleared out entirely.
foreach ($idx in 1,2) {
    $ws = $wb.Worksheets.Item($idx)

    $ws.Range("B2").Value = "FreeText"
    $ws.Range("A2").Value = "fdfuu42a22321c123a8_test"
    $ws.Range("C2").Value = "Test question fdfuu42a22321c123a8_test"
    $ws.Range("D2").Value = "Test question fdfuu42a22321c123a8_test"

    $ws.Range("A3").Value = "fdfzz42a66321c123a8_test"
    $ws.Range("B3").Value = "FreeText"
    $ws.Range("C3").Value = "Test question fdfzz42a66321c123a8_test"
    $ws.Range("D3").Value = "Test question fdfzz42a66321c123a8_test"

    $ws.Range("A4:D4").ClearContents()
}

# Rename the sheets (case/number cleanup)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "test_new_survey_import 2"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "test_new_survey_import 3"

# Restore selections on each sheet / which tab is active
$ws1.Activate()
$ws1.Range("A2:D3").Select()

$ws2.Activate()
$ws2.Range("F18").Select()
